$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("researchMeasures")

# --- Update the notes_diet_mood_etc text in Z34 (append new sentences) ---
$ws.Range("Z34").Value = 'Woke up to have a BM around 1130 pm after an hour of sleep, it was a solid constipated small BM, then went back to bed and woke up at 5 am approximately and laid in bed until 530 am, got up made the babies their food, my coffee, and had the rest of last night''s BM also a solid reg size one. Together a lg BM, might be the Aldi''s Happy Farm Cheese making it like that consistency, but also I haven''t been drinking much water, only while working out, and that was only 1 bottle as it was cold during the workout and after. Took my measurements at 620 am. I tried the batteries I bought yesterday the AAA ones for my calculater, energizer brand but they make the screen dark and worse than the other batteries that are generic. So I put back in the generic ones. They cost me $10 just for those batteries too. I will see if changing the 3v battery will fix this problem but I need to buy that one too. Shortly after 6:45 am had a reg BM that feels like it cleared up the other solid waste hanging around the rectum with the normal vegetarian consistency/texture of my reg BMs. My courses in genetics and general chemistry AKA organic chemistry start this Thursday. I have been reading ahead and plan to do more today before lifting weights. No cardio today, I think I am just going to switch or alternate between cardio and weight lifting days. I have been so far, during the middle of this research. No menstruation yet but I usually get it before 28 days. Clearly, not sexually active and haven''t been for years as digestion is only body function interested in for body conditioning and maintanence. So far no obvious changes in fibroid waistline changes, but all working out and using the waist trimmer have helped keep it compacted and small as I can look down and see my pubes. I couldn''t about a few months ago, still not where they start but getting close. At least I am healthy and strong. Ate break fast a little after 7 am, with 2 eggs scrambled in 2 tbsp sourcream and cooked in 2 tbsp olive oil, and 2 corn tortilla quesadillas with the Guerrero brand tortillas and the Winco low skim mozzarella cheese. And an orange around 730 am. Read my chemistry notes.'

# --- Add the new dailyFoodConsumed note in AA34 ---
$ws.Range("AA34").Value = "2 eggs scrambled with `n(140`t10`t3`t12`t0`t0`t140)`nolive oil 2 tbsp and `n(120`t14`t2`t0`t0`t0`t0)`nsourcream 2 tbsp Daisy brand last of it`n(60`t5`t3.5`t1`t0`t1`t15)`n2 corn tortilla (Guerrero Brand) with Winco low skim mozzarella shredded cheese and paprika`n4 corn tortillas (200`t  2`t0`t4`t42`t4`t40)`n1/2 cup mozzarella cheese (160`t 10`t7`t12`t2`t0`t380)`n1 orange (81`t0`t0`t2`t21`t4`t2)`n"

# --- Macro nutrient totals for the day (AB34:AH34) ---
$ws.Range("AB34").Formula = "=140+120+60+200+160+81"
$ws.Range("AC34").Formula = "=10+14+5+2+10+0"
$ws.Range("AD34").Formula = "=3+2+3.5+0+7+0"
$ws.Range("AE34").Formula = "=12+0+1+4+12+2"
$ws.Range("AF34").Formula = "=0+0+0+42+2+21"
$ws.Range("AG34").Formula = "=0+0+1+4+0+4"
$ws.Range("AH34").Formula = "=140+0+15+40+380+2"

# --- Ratios relative to total calories (AI34:AN34) ---
$ws.Range("AI34").Formula = '=$AC34/$AB34'
$ws.Range("AJ34").Formula = '=$AD34/$AB34'
$ws.Range("AK34").Formula = '=$AE34/$AB34'
$ws.Range("AL34").Formula = '=$AF34/$AB34'
$ws.Range("AM34").Formula = '=$AG34/$AB34'
$ws.Range("AN34").Formula = '=$AH34/$AB34'

# --- coffee_cups / BM / Menstruation ---
$ws.Range("AO34").Value = 3
$ws.Range("AP34").Value = 3
$ws.Range("AQ34").Value = 0

# --- weight lifting increase/decrease flags and amounts, waist trimmer, etc (AS34:BD34) ---
$ws.Range("AS34").Value = 0
$ws.Range("AT34").Value = 0
$ws.Range("AU34").Value = 0
$ws.Range("AV34").Value = 0
$ws.Range("AW34").Value = 31
$ws.Range("AX34").Value = 1
$ws.Range("AY34").Value = 7
$ws.Range("AZ34").Value = 1
$ws.Range("BA34").Value = 1
$ws.Range("BB34").Value = 1
$ws.Range("BC34").Value = 1
$ws.Range("BD34").Value = 1

# Match styling (number format + alignment) for the cells whose style differs
# from the worksheet column default so the underlying style index lines up
# with style "7" (0 decimal, horizontal/vertical centered, fill applied).
$styledRange = $ws.Range("AS34:BD34")
$styledRange.NumberFormat = "0"
$styledRange.HorizontalAlignment = -4108
$styledRange.VerticalAlignment = -4160

# Re-assert the original row height; adding the long note text can trigger
# automatic row auto-fit which would otherwise change the row's height.
$ws.Rows.Item(34).RowHeight = 20.1

# --- Update the NutritionalData sheet selection/view ---
$ws3 = $wb.Worksheets.Item("NutritionalData")
$ws3.Range("B7:H7").Select()

# --- Restore the researchMeasures sheet as active and set the new selection ---
$ws.Activate()
$ws.Range("AN35").Select()
